$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "67.632.01"
Set-TextValue $ws.Range("E2") "  -1.37%  "

Set-TextValue $ws.Range("D3") "3.770.75"
Set-TextValue $ws.Range("E3") "  -2.15%  "

Set-TextValue $ws.Range("D5") "596.40"
Set-TextValue $ws.Range("E5") "  -1.08%  "

Set-TextValue $ws.Range("D6") "168.90"
Set-TextValue $ws.Range("E6") "  -0.31%  "

Set-TextValue $ws.Range("D7") "3.769.25"
Set-TextValue $ws.Range("E7") "  -2.17%  "

Set-TextValue $ws.Range("E8") "  -0.01%  "

Set-TextValue $ws.Range("E9") "  -0.90%  "

Set-TextValue $ws.Range("D10") "0.164"
Set-TextValue $ws.Range("E10") "  -1.72%  "

Set-TextValue $ws.Range("D11") "6.46"
Set-TextValue $ws.Range("E11") "  -0.62%  "

Set-TextValue $ws.Range("D12") "0.452"
Set-TextValue $ws.Range("E12") "  -1.20%  "

Set-TextValue $ws.Range("D13") "0.0000275"
Set-TextValue $ws.Range("E13") "  +1.59%  "

Set-TextValue $ws.Range("D14") "36.50"
Set-TextValue $ws.Range("E14") "  -1.78%  "

Set-TextValue $ws.Range("D15") "4.404.49"
Set-TextValue $ws.Range("E15") "  -2.12%  "

Set-TextValue $ws.Range("D16") "3.772.76"
Set-TextValue $ws.Range("E16") "  -1.79%  "

Set-TextValue $ws.Range("D17") "18.57"
Set-TextValue $ws.Range("E17") "  +0.44%  "

Set-TextValue $ws.Range("D18") "67.641.42"
Set-TextValue $ws.Range("E18") "  -1.52%  "

Set-TextValue $ws.Range("D19") "7.18"
Set-TextValue $ws.Range("E19") "  -3.06%  "

Set-TextValue $ws.Range("E20") "  +0.70%  "

Set-TextValue $ws.Range("D21") "10.57"
Set-TextValue $ws.Range("E21") "  -4.85%  "

Set-TextValue $ws.Range("D22") "467.52"
Set-TextValue $ws.Range("E22") "  -0.78%  "

Set-TextValue $ws.Range("D23") "0.718"
Set-TextValue $ws.Range("E23") "  -2.41%  "

Set-TextValue $ws.Range("D24") "0.0000147"
Set-TextValue $ws.Range("E24") "  -9.81%  "

Set-TextValue $ws.Range("D25") "83.70"
Set-TextValue $ws.Range("E25") "  +0.22%  "

Set-TextValue $ws.Range("D26") "2.20"
Set-TextValue $ws.Range("E26") "  -1.42%  "

Set-TextValue $ws.Range("E27") "  +0.00%  "

Set-TextValue $ws.Range("D28") "10.26"
Set-TextValue $ws.Range("E28") "  +1.01%  "

Set-TextValue $ws.Range("E29") "  -0.15%  "

Set-TextValue $ws.Range("E30") "  -2.06%  "

Set-TextValue $ws.Range("D31") "3.923.43"
Set-TextValue $ws.Range("E31") "  -2.01%  "

Set-TextValue $ws.Range("D32") "7.62"
Set-TextValue $ws.Range("E32") "  -1.20%  "

Set-TextValue $ws.Range("D33") "30.46"
Set-TextValue $ws.Range("E33") "  -3.53%  "

Set-TextValue $ws.Range("D34") "2.22"
Set-TextValue $ws.Range("E34") "  -3.93%  "

Set-TextValue $ws.Range("D35") "9.11"
Set-TextValue $ws.Range("E35") "  -2.91%  "

Set-TextValue $ws.Range("D36") "3.732.71"
Set-TextValue $ws.Range("E36") "  -2.25%  "

Set-TextValue $ws.Range("D37") "3.83"
Set-TextValue $ws.Range("E37") "  +2.70%  "

Set-TextValue $ws.Range("E38") "  -0.94%  "

Set-TextValue $ws.Range("E39") "  -1.05%  "

Set-TextValue $ws.Range("E40") "  -1.99%  "

Set-TextValue $ws.Range("D41") "5.78"
Set-TextValue $ws.Range("E41") "  -2.78%  "

Set-TextValue $ws.Range("E42") "  +0.05%  "

Set-TextValue $ws.Range("D43") "0.311"
Set-TextValue $ws.Range("E43") "  -1.47%  "

Set-TextValue $ws.Range("E44") "  +0.01%  "

Set-TextValue $ws.Range("D45") "8.67"
Set-TextValue $ws.Range("E45") "  -0.53%  "

Set-TextValue $ws.Range("E46") "  -2.14%  "

Set-TextValue $ws.Range("D47") "45.76"
Set-TextValue $ws.Range("E47") "  -2.87%  "

Set-TextValue $ws.Range("D48") "394.87"
Set-TextValue $ws.Range("E48") "  -5.29%  "

Set-TextValue $ws.Range("D49") "0.000270"
Set-TextValue $ws.Range("E49") "  -8.70%  "

Set-TextValue $ws.Range("D50") "140.21"
Set-TextValue $ws.Range("E50") "  -0.88%  "

Set-TextValue $ws.Range("D51") "0.0352"
Set-TextValue $ws.Range("E51") "  -2.43%  "
